# Remove the 'culture_collection' MIGS attribute column (column Y) from the
# template. Deleting the column shifts the worksheet data (headers + shared
# strings) correctly, but this runtime does not move cell comments when a
# column is deleted, so the comments anchored on row 15 (the header/comment
# row) have to be re-pointed by hand to match their shifted columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellComment {
    param(
        $Sheet,
        [string]$Ref,
        $CommentText
    )
    $cell = $Sheet.Range($Ref)
    if ($cell.Comment -ne $null) {
        $cell.Comment.Delete()
    }
    if ($CommentText -ne $null) {
        $cell.AddComment($CommentText) | Out-Null
    }
}

# Delete column Y ('culture_collection'); everything to its right shifts one
# column to the left (sheet cell values + shared strings are handled by
# Excel automatically).
$ws.Columns("Y").Delete()

# Re-point the row-15 comments: each surviving column now holds the comment
# that used to belong to the column immediately to its right (before the
# delete). The final comment (old BT15 / 'trophic_level') is removed below
# once its text has been moved into BS15.
Set-CellComment $ws 'Y15' 'history of dermatology disorders; can include multiple disorders'
Set-CellComment $ws 'Z15' 'dominant hand of the subject'
Set-CellComment $ws 'AA15' 'ethnicity of the subject'
Set-CellComment $ws 'AB15' 'Plasmids that have significance phenotypic consequence'
Set-CellComment $ws 'AC15' 'Health or disease status of sample at time of collection'
Set-CellComment $ws 'AD15' 'Age of host at the time of sampling'
Set-CellComment $ws 'AE15' 'body mass index of the host, calculated as weight/(height)squared'
Set-CellComment $ws 'AF15' 'substance produced by the host, e.g. stool, mucus, where the sample was obtained from'
Set-CellComment $ws 'AG15' 'core body temperature of the host when sample was collected'
Set-CellComment $ws 'AH15' 'type of diet depending on the sample for animals omnivore, herbivore etc., for humans high-fat, meditteranean etc.; can include multiple diet types'
Set-CellComment $ws 'AI15' 'Name of relevant disease, e.g. Salmonella gastroenteritis. For the controlled vocabulary, please see Human Disease Ontology, http://bioportal.bioontology.org/ontologies/1009 or MeSH, http://www.ncbi.nlm.nih.gov/mesh'
Set-CellComment $ws 'AJ15' $null
Set-CellComment $ws 'AK15' $null
Set-CellComment $ws 'AL15' 'the height of subject'
Set-CellComment $ws 'AM15' 'content of last meal and time since feeding; can include multiple values'
Set-CellComment $ws 'AN15' 'most frequent job performed by subject'
Set-CellComment $ws 'AO15' $null
Set-CellComment $ws 'AP15' 'resting pulse of the host, measured as beats per minute'
Set-CellComment $ws 'AQ15' 'Gender or physical sex of the host'
Set-CellComment $ws 'AR15' 'a unique identifier by which each subject can be referred to, de-identified, e.g. #131'
Set-CellComment $ws 'AS15' 'NCBI taxonomy ID of the host, e.g. 9606'
Set-CellComment $ws 'AT15' 'Type of tissue the initial sample was taken from. Controlled vocabulary, http://bioportal.bioontology.org/ontologies/1005'
Set-CellComment $ws 'AU15' 'total mass of the host at collection, the unit depends on host'
Set-CellComment $ws 'AV15' 'can include multiple medication codes'
Set-CellComment $ws 'AW15' 'Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.'
Set-CellComment $ws 'AX15' 'A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html'
Set-CellComment $ws 'AY15' 'whether full medical history was collected'
Set-CellComment $ws 'AZ15' 'any other measurement performed or parameter collected, that is not listed here'
Set-CellComment $ws 'BA15' 'total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts'
Set-CellComment $ws 'BB15' 'oxygenation status of sample'
Set-CellComment $ws 'BC15' 'To what is the entity pathogenic'
Set-CellComment $ws 'BD15' 'type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types'
Set-CellComment $ws 'BE15' 'Primary publication or genome report in the form of pubmed ID, DOI or URL'
Set-CellComment $ws 'BF15' 'Method or device employed for collecting sample'
Set-CellComment $ws 'BG15' 'Processing applied to the sample during or after isolation'
Set-CellComment $ws 'BH15' 'salinity of sample, i.e. measure of total salt concentration'
Set-CellComment $ws 'BI15' 'Amount or size of sample (volume, mass or area) that was collected'
Set-CellComment $ws 'BJ15' 'duration for which sample was stored'
Set-CellComment $ws 'BK15' 'location at which sample was stored, usually name of a specific freezer/room'
Set-CellComment $ws 'BL15' 'temperature at which sample was stored, e.g. -80'
Set-CellComment $ws 'BM15' 'volume (mL) or weight (g) of sample processed for DNA extraction'
Set-CellComment $ws 'BN15' 'unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.'
Set-CellComment $ws 'BO15' 'Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier'
Set-CellComment $ws 'BP15' 'Information about the genetic distinctness of the lineage (eg., biovar, serovar)'
Set-CellComment $ws 'BQ15' 'temperature of the sample at time of sampling'
Set-CellComment $ws 'BR15' 'specification of the time since last wash'
Set-CellComment $ws 'BS15' 'Feeding position in food chain (eg., chemolithotroph)'

# The old last column (BT, 'trophic_level') has no column to its right to
# inherit from anymore now that the data only spans to BS, so drop it.
Set-CellComment $ws 'BT15' $null

Write-Output "culture_collection column removed; row 15 comments re-pointed."
